# schema-repo-metadata.xlsx — reorder git metadata rows to url/branch/revision
# and register the new (blank) trailing row; keep the data-validation rules
# attached to the cell that now holds the matching label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schema repo metadata")

# --- Reorder the metadata rows: Url, Branch, Revision (was Branch, Revision, Url) ---
$ws.Cells.Item(1, 1).Value = "Url"
$ws.Cells.Item(1, 2).Value = "https://github.com/KarrLab/test_repo"
$ws.Cells.Item(2, 1).Value = "Branch"
$ws.Cells.Item(2, 2).Value = "master"
$ws.Cells.Item(3, 1).Value = "Revision"
$ws.Cells.Item(3, 2).Value = "feb697e6e90f06da7a0aecb873ad8e235053a97f"

# --- Register the new trailing blank row (dimension grows from A1:B3 to A1:B4) ---
$ws.Rows.Item(4).RowHeight = 15
$ws.Cells.Item(4, 1).Style = "Normal"

# --- Keep each data-validation rule's title pinned to its label, which now
#     lives on a different row: Url -> B1, Branch -> B2, Revision -> B3 ---
$ws.Range("B1").Validation.ErrorTitle = "Url"
$ws.Range("B1").Validation.InputTitle = "Url"
$ws.Range("B2").Validation.ErrorTitle = "Branch"
$ws.Range("B2").Validation.InputTitle = "Branch"
$ws.Range("B3").Validation.ErrorTitle = "Revision"
$ws.Range("B3").Validation.InputTitle = "Revision"

# --- Move the active selection from B3 to A2 ---
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null

# --- Best-effort: move the workbook window's on-screen position ---
$win = $excel.ActiveWindow
$win.Left = 5440
$win.Top = 8860
